$d = $word.ActiveDocument

$replacements = @(
    @("62÷3=", "77÷4="),
    @("99÷3=", "77÷7="),
    @("24÷6=", "60÷9="),
    @("97÷5=", "33÷3="),
    @("21÷8=", "50÷5="),
    @("53÷8=", "57÷6="),
    @("32÷7=", "13÷6="),
    @("11÷6=", "18÷6="),
    @("79÷6=", "16÷6="),
    @("61÷3=", "32÷2="),
    @("36÷6=", "57÷3="),
    @("84÷9=", "14÷2="),
    @("32÷6=", "98÷7="),
    @("69÷6=", "37÷4="),
    @("13÷3=", "65÷3="),
    @("59÷5=", "76÷9="),
    @("10÷8=", "61÷2="),
    @("60÷6=", "79÷6="),
    @("64÷6=", "39÷9="),
    @("43÷8=", "67÷6="),
    @("81÷4=", "44÷7="),
    @("77÷8=", "19÷9="),
    @("70÷7=", "67÷2="),
    @("51÷6=", "42÷9="),
    @("84÷3=", "22÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
